$d = $word.ActiveDocument

# 1. Fix "Team menbers" -> "Team members"
$d.Content.Find.Execute("Team menbers:", $false, $false, $false, $false, $false, $true, 1, $false, "Team members:", 2)

# 2. Fix "aimming" -> "aiming"
$d.Content.Find.Execute("aimming", $false, $false, $false, $false, $false, $true, 1, $false, "aiming", 2)

# 3. "help themselves be" -> "help them be"
$d.Content.Find.Execute("help themselves be", $false, $false, $false, $false, $false, $true, 1, $false, "help them be", 2)

# 4. "completed(such" -> "completed (such"
$d.Content.Find.Execute("completed(such", $false, $false, $false, $false, $false, $true, 1, $false, "completed (such", 2)

# 5. "can be moviated" -> "can be motivated"
$d.Content.Find.Execute("can be moviated", $false, $false, $false, $false, $false, $true, 1, $false, "can be motivated", 2)

# 6. "our lfe-time" -> "our life-time"
$d.Content.Find.Execute("our lfe-time", $false, $false, $false, $false, $false, $true, 1, $false, "our life-time", 2)

# 7. "studying paterners" -> "studying partners"
$d.Content.Find.Execute("studying paterners", $false, $false, $false, $false, $false, $true, 1, $false, "studying partners", 2)

# 8. Tidy up the run break in the middle of the Purpose paragraph (no text change)
$d.Content.Find.Execute("to solve the problem that people", $false, $false, $false, $false, $false, $true, 1, $false, "to solve the problem that people", 2)

# 9. Tidy up the run break in the Marketing paragraph (no text change)
$d.Content.Find.Execute("have interests in sharing rooms", $false, $false, $false, $false, $false, $true, 1, $false, "have interests in sharing rooms", 2)

# 10. Tidy up the run break after "Technical Requirements" heading (no text change)
$d.Content.Find.Execute("Technical Requirements:", $false, $false, $false, $false, $false, $true, 1, $false, "Technical Requirements:", 2)

# 11. "Hypervizor:" -> "Hypervizor: VirtualBox"
$d.Content.Find.Execute("Hypervizor:", $false, $false, $false, $false, $false, $true, 1, $false, "Hypervizor: VirtualBox", 2)

# 12. "Support Libraries: Mcrypt, others" -> "Support Libraries: Mcrypt, Sentinel, others"
$d.Content.Find.Execute("Mcrypt, others", $false, $false, $false, $false, $false, $true, 1, $false, "Mcrypt, Sentinel, others", 2)

# 13. "Access to server: Git" -> "Access to server: Git, SSH"
$d.Content.Find.Execute("Access to server: Git", $false, $false, $false, $false, $false, $true, 1, $false, "Access to server: Git, SSH", 2)

# 14. Security / Authentication text replacement
$d.Content.Find.Execute("Security / Authentication: How do you plan to encrypt passwords, and authenticate users?", $false, $false, $false, $false, $false, $true, 1, $false, "Security / Authentication: Using 3rd party library (Sentinel)", 2)

# 14b. Make the "rd" in "3rd" superscript
$r = $d.Content
$r.Find.Execute("Using 3rd")
$r.Start = $r.End - 2
$r.Font.Superscript = $true
